$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @(2,7,[double]"25.112885"),
  @(2,8,[double]"75.33865499999999"),
  @(2,9,[double]"0.4760809326504767"),
  @(2,10,[double]"0.4760809326504767"),
  @(2,11,[double]"2"),
  @(2,12,[double]"0.6666666666666666"),
  @(2,13,[double]"0.2262196666666666"),
  @(2,14,[double]"0.6786589999999999"),
  @(2,15,[double]"0.03145179203784564"),
  @(2,16,[double]"0.03145179203784564"),
  @(2,17,[double]"5.681028473738331"),
  @(2,18,[double]"51.12925626364498"),
  @(2,19,[double]"0.01497359848690639"),
  @(2,20,[double]"0.01497359848690639"),
  @(3,7,[double]"25.112885"),
  @(3,8,[double]"75.33865499999999"),
  @(3,9,[double]"0.4760809326504767"),
  @(3,10,[double]"0.4760809326504767"),
  @(3,15,[double]"0.9636438974901603"),
  @(3,16,[double]"0.9636438974901604"),
  @(3,17,[double]"174.05966609465"),
  @(3,18,[double]"1566.53699485185"),
  @(3,19,[double]"0.4587724854600559"),
  @(3,20,[double]"0.4587724854600559"),
  @(4,7,[double]"25.112885"),
  @(4,8,[double]"75.33865499999999"),
  @(4,9,[double]"0.4760809326504767"),
  @(4,10,[double]"0.4760809326504767"),
  @(4,11,[double]"1"),
  @(4,12,[double]"0.3333333333333333"),
  @(4,13,[double]"0.03527466666666667"),
  @(4,14,[double]"0.105824"),
  @(4,15,[double]"0.004904310471994002"),
  @(4,16,[double]"0.004904310471994003"),
  @(4,17,[double]"0.8858486474133332"),
  @(4,18,[double]"7.972637826719999"),
  @(4,19,[double]"0.002334848703514404"),
  @(4,20,[double]"0.002334848703514405"),
  @(5,9,[double]"0.0009058021671219375"),
  @(5,10,[double]"0.0009058021671219374"),
  @(5,11,[double]"2"),
  @(5,12,[double]"0.6666666666666666"),
  @(5,13,[double]"0.2262196666666666"),
  @(5,14,[double]"0.6786589999999999"),
  @(5,15,[double]"0.03145179203784564"),
  @(5,16,[double]"0.03145179203784564"),
  @(5,17,[double]"0.01080885107988889"),
  @(5,18,[double]"0.09727965971899999"),
  @(5,19,[double]"2.848910138774908E-05"),
  @(5,20,[double]"2.848910138774908E-05"),
  @(6,9,[double]"0.0009058021671219375"),
  @(6,10,[double]"0.0009058021671219374"),
  @(6,15,[double]"0.9636438974901603"),
  @(6,16,[double]"0.9636438974901604"),
  @(6,19,[double]"0.0008728707306804174"),
  @(6,20,[double]"0.0008728707306804174"),
  @(7,9,[double]"0.0009058021671219375"),
  @(7,10,[double]"0.0009058021671219374"),
  @(7,11,[double]"1"),
  @(7,12,[double]"0.3333333333333333"),
  @(7,13,[double]"0.03527466666666667"),
  @(7,14,[double]"0.105824"),
  @(7,15,[double]"0.004904310471994002"),
  @(7,16,[double]"0.004904310471994003"),
  @(7,17,[double]"0.001685435331555556"),
  @(7,18,[double]"0.015168917984"),
  @(7,19,[double]"4.442335053770979E-06"),
  @(7,20,[double]"4.44233505377098E-06"),
  @(8,7,[double]"1.993125666666667"),
  @(8,8,[double]"5.979377"),
  @(8,9,[double]"0.03778495088382995"),
  @(8,10,[double]"0.03778495088382995"),
  @(8,11,[double]"2"),
  @(8,12,[double]"0.6666666666666666"),
  @(8,13,[double]"0.2262196666666666"),
  @(8,14,[double]"0.6786589999999999"),
  @(8,15,[double]"0.03145179203784564"),
  @(8,16,[double]"0.03145179203784564"),
  @(8,17,[double]"0.4508842239381111"),
  @(8,18,[double]"4.057958015443"),
  @(8,19,[double]"0.001188404417358432"),
  @(8,20,[double]"0.001188404417358432"),
  @(9,7,[double]"1.993125666666667"),
  @(9,8,[double]"5.979377"),
  @(9,9,[double]"0.03778495088382995"),
  @(9,10,[double]"0.03778495088382995"),
  @(9,15,[double]"0.9636438974901603"),
  @(9,16,[double]"0.9636438974901604"),
  @(9,17,[double]"13.81453337697667"),
  @(9,18,[double]"124.33080039279"),
  @(9,19,[double]"0.03641123733616818"),
  @(9,20,[double]"0.03641123733616818"),
  @(10,7,[double]"1.993125666666667"),
  @(10,8,[double]"5.979377"),
  @(10,9,[double]"0.03778495088382995"),
  @(10,10,[double]"0.03778495088382995"),
  @(10,11,[double]"1"),
  @(10,12,[double]"0.3333333333333333"),
  @(10,13,[double]"0.03527466666666667"),
  @(10,14,[double]"0.105824"),
  @(10,15,[double]"0.004904310471994002"),
  @(10,16,[double]"0.004904310471994003"),
  @(10,17,[double]"0.07030684351644445"),
  @(10,18,[double]"0.6327615916480001"),
  @(10,19,[double]"0.0001853091303033463"),
  @(10,20,[double]"0.0001853091303033463"),
  @(11,7,[double]"0.147857"),
  @(11,8,[double]"0.443571"),
  @(11,9,[double]"0.002803019185525739"),
  @(11,10,[double]"0.002803019185525739"),
  @(11,11,[double]"2"),
  @(11,12,[double]"0.6666666666666666"),
  @(11,13,[double]"0.2262196666666666"),
  @(11,14,[double]"0.6786589999999999"),
  @(11,15,[double]"0.03145179203784564"),
  @(11,16,[double]"0.03145179203784564"),
  @(11,17,[double]"0.03344816125433333"),
  @(11,18,[double]"0.301033451289"),
  @(11,19,[double]"8.815997650124702E-05"),
  @(11,20,[double]"8.815997650124702E-05"),
  @(12,7,[double]"0.147857"),
  @(12,8,[double]"0.443571"),
  @(12,9,[double]"0.002803019185525739"),
  @(12,10,[double]"0.002803019185525739"),
  @(12,15,[double]"0.9636438974901603"),
  @(12,16,[double]"0.9636438974901604"),
  @(12,17,[double]"1.02481017413"),
  @(12,18,[double]"9.223291567170001"),
  @(12,19,[double]"0.002701112332679718"),
  @(12,20,[double]"0.002701112332679718"),
  @(13,7,[double]"0.147857"),
  @(13,8,[double]"0.443571"),
  @(13,9,[double]"0.002803019185525739"),
  @(13,10,[double]"0.002803019185525739"),
  @(13,11,[double]"1"),
  @(13,12,[double]"0.3333333333333333"),
  @(13,13,[double]"0.03527466666666667"),
  @(13,14,[double]"0.105824"),
  @(13,15,[double]"0.004904310471994002"),
  @(13,16,[double]"0.004904310471994003"),
  @(13,17,[double]"0.005215606389333334"),
  @(13,18,[double]"0.04694045750400001"),
  @(13,19,[double]"1.374687634477398E-05"),
  @(13,20,[double]"1.374687634477399E-05"),
  @(14,7,[double]"22.78687466666667"),
  @(14,8,[double]"68.360624"),
  @(14,9,[double]"0.4319852754271836"),
  @(14,10,[double]"0.4319852754271836"),
  @(14,11,[double]"2"),
  @(14,12,[double]"0.6666666666666666"),
  @(14,13,[double]"0.2262196666666666"),
  @(14,14,[double]"0.6786589999999999"),
  @(14,15,[double]"0.03145179203784564"),
  @(14,16,[double]"0.03145179203784564"),
  @(14,17,[double]"5.154839191468444"),
  @(14,18,[double]"46.393552723216"),
  @(14,19,[double]"0.01358671104614725"),
  @(14,20,[double]"0.01358671104614725"),
  @(15,7,[double]"22.78687466666667"),
  @(15,8,[double]"68.360624"),
  @(15,9,[double]"0.4319852754271836"),
  @(15,10,[double]"0.4319852754271836"),
  @(15,15,[double]"0.9636438974901603"),
  @(15,16,[double]"0.9636438974901604"),
  @(15,17,[double]"157.9378791333867"),
  @(15,18,[double]"1421.44091220048"),
  @(15,19,[double]"0.4162799744710116"),
  @(15,20,[double]"0.4162799744710117"),
  @(16,7,[double]"22.78687466666667"),
  @(16,8,[double]"68.360624"),
  @(16,9,[double]"0.4319852754271836"),
  @(16,10,[double]"0.4319852754271836"),
  @(16,11,[double]"1"),
  @(16,12,[double]"0.3333333333333333"),
  @(16,13,[double]"0.03527466666666667"),
  @(16,14,[double]"0.105824"),
  @(16,15,[double]"0.004904310471994002"),
  @(16,16,[double]"0.004904310471994003"),
  @(16,17,[double]"0.8037994082417779"),
  @(16,18,[double]"7.234194674176"),
  @(16,19,[double]"0.00211858991002475"),
  @(16,20,[double]"0.002118589910024751"),
  @(17,7,[double]"2.660670333333333"),
  @(17,8,[double]"7.982011"),
  @(17,9,[double]"0.05044001968586199"),
  @(17,10,[double]"0.05044001968586199"),
  @(17,11,[double]"2"),
  @(17,12,[double]"0.6666666666666666"),
  @(17,13,[double]"0.2262196666666666"),
  @(17,14,[double]"0.6786589999999999"),
  @(17,15,[double]"0.03145179203784564"),
  @(17,16,[double]"0.03145179203784564"),
  @(17,17,[double]"0.6018959559165554"),
  @(17,18,[double]"5.417063603248999"),
  @(17,19,[double]"0.001586429009544572"),
  @(17,20,[double]"0.001586429009544572"),
  @(18,7,[double]"2.660670333333333"),
  @(18,8,[double]"7.982011"),
  @(18,9,[double]"0.05044001968586199"),
  @(18,10,[double]"0.05044001968586199"),
  @(18,15,[double]"0.9636438974901603"),
  @(18,16,[double]"0.9636438974901604"),
  @(18,17,[double]"18.44134554066333"),
  @(18,18,[double]"165.97210986597"),
  @(18,19,[double]"0.04860621715956447"),
  @(18,20,[double]"0.04860621715956447"),
  @(19,7,[double]"2.660670333333333"),
  @(19,8,[double]"7.982011"),
  @(19,9,[double]"0.05044001968586199"),
  @(19,10,[double]"0.05044001968586199"),
  @(19,11,[double]"1"),
  @(19,12,[double]"0.3333333333333333"),
  @(19,13,[double]"0.03527466666666667"),
  @(19,14,[double]"0.105824"),
  @(19,15,[double]"0.004904310471994002"),
  @(19,16,[double]"0.004904310471994003"),
  @(19,17,[double]"0.09385425911822223"),
  @(19,18,[double]"0.844688332064"),
  @(19,19,[double]"0.0002473735167529566"),
  @(19,20,[double]"0.0002473735167529566")
)

foreach ($u in $updates) {
  $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}